$wb = $excel.ActiveWorkbook

function Set-CellValue($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

function Clear-CellValue($ws, $row, $col) {
    $ws.Cells.Item($row, $col).ClearContents()
}

$ws = $wb.Worksheets.Item("ALC")
# Row 70
Set-CellValue $ws 70 8 1621.0834
Set-CellValue $ws 70 9 2325.5
Set-CellValue $ws 70 10 916.6667
Set-CellValue $ws 70 11 6976.5
Set-CellValue $ws 70 12 2750.0001
Set-CellValue $ws 70 13 -6706.5
Set-CellValue $ws 70 14 -3290.0001

# Row 73
Set-CellValue $ws 73 8 1621.0834
Set-CellValue $ws 73 9 2325.5
Set-CellValue $ws 73 10 916.6667
Set-CellValue $ws 73 11 6976.5
Set-CellValue $ws 73 12 2750.0001
Set-CellValue $ws 73 13 -6040.5
Set-CellValue $ws 73 14 -4622.0001

# Row 74
Set-CellValue $ws 74 8 2721.2
Set-CellValue $ws 74 9 2298.2222
Set-CellValue $ws 74 11 2298.2222
Set-CellValue $ws 74 13 -1362.2222

# Row 77
Set-CellValue $ws 77 8 2721.2
Set-CellValue $ws 77 9 2298.2222
Set-CellValue $ws 77 11 11491.111
Set-CellValue $ws 77 13 -6811.111000000001

# Row 132
Set-CellValue $ws 132 8 28851972
Set-CellValue $ws 132 9 31921250
Set-CellValue $ws 132 11 95763750
Set-CellValue $ws 132 13 -95761220

# Row 138
Set-CellValue $ws 138 8 3673.7334
Set-CellValue $ws 138 9 1287.7906
Set-CellValue $ws 138 10 6879.8438
Set-CellValue $ws 138 11 3863.3718
Set-CellValue $ws 138 12 20639.5314
Set-CellValue $ws 138 14 -30919.5314
Set-CellValue $ws 138 13 1276.6282

# Row 141
Set-CellValue $ws 141 8 5474.4
Set-CellValue $ws 141 9 1806.5
Set-CellValue $ws 141 11 5419.5
Set-CellValue $ws 141 13 -239.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
Set-CellValue $ws 32 8 6295595.5
Set-CellValue $ws 32 9 5457.5654
Set-CellValue $ws 32 10 47630788
Set-CellValue $ws 32 11 5457.5654
Set-CellValue $ws 32 12 47630788
Set-CellValue $ws 32 13 -5170.5654
Set-CellValue $ws 32 14 -47631362

# Row 45
Set-CellValue $ws 45 8 12860058
Set-CellValue $ws 45 9 19667532
Set-CellValue $ws 45 10 1494.3334
Set-CellValue $ws 45 11 19667532
Set-CellValue $ws 45 12 1494.3334
Set-CellValue $ws 45 13 -19667155
Set-CellValue $ws 45 14 -2248.3334

# Row 102
Set-CellValue $ws 102 8 0
Set-CellValue $ws 102 9 0
Set-CellValue $ws 102 10 0
Set-CellValue $ws 102 11 0
Set-CellValue $ws 102 12 0
Clear-CellValue $ws 102 13
Clear-CellValue $ws 102 14

$ws = $wb.Worksheets.Item("BSM")
# Row 20
Set-CellValue $ws 20 8 1092.68
Set-CellValue $ws 20 9 1234.5714
Set-CellValue $ws 20 10 912.0909
Set-CellValue $ws 20 11 1234.5714
Set-CellValue $ws 20 12 912.0909
Set-CellValue $ws 20 13 -987.5714
Set-CellValue $ws 20 14 -1406.0909

# Row 59
Set-CellValue $ws 59 8 49290
Set-CellValue $ws 59 10 49290
Set-CellValue $ws 59 12 49290
Set-CellValue $ws 59 14 -50984

# Row 86
Set-CellValue $ws 86 8 1389.303
Set-CellValue $ws 86 9 1416.4073
Set-CellValue $ws 86 10 1267.3334
Set-CellValue $ws 86 11 1416.4073
Set-CellValue $ws 86 12 1267.3334
Set-CellValue $ws 86 13 -293.4073000000001
Set-CellValue $ws 86 14 -3513.3334

# Row 89
Set-CellValue $ws 89 8 1389.303
Set-CellValue $ws 89 9 1416.4073
Set-CellValue $ws 89 10 1267.3334
Set-CellValue $ws 89 11 7082.0365
Set-CellValue $ws 89 12 6336.666999999999
Set-CellValue $ws 89 13 -1466.0365
Set-CellValue $ws 89 14 -17568.667

# Row 99
Set-CellValue $ws 99 8 1325.102
Set-CellValue $ws 99 9 1199
Set-CellValue $ws 99 10 1524.2106
Set-CellValue $ws 99 11 1199
Set-CellValue $ws 99 12 1524.2106
Set-CellValue $ws 99 13 299
Set-CellValue $ws 99 14 -4520.2106

$ws = $wb.Worksheets.Item("CRP")
# Row 31
Set-CellValue $ws 31 8 1150.8103
Set-CellValue $ws 31 9 815.3913
Set-CellValue $ws 31 10 1371.2285
Set-CellValue $ws 31 11 815.3913
Set-CellValue $ws 31 12 1371.2285
Set-CellValue $ws 31 13 -520.3913
Set-CellValue $ws 31 14 -1961.2285

# Row 34
Set-CellValue $ws 34 8 1150.8103
Set-CellValue $ws 34 9 815.3913
Set-CellValue $ws 34 10 1371.2285
Set-CellValue $ws 34 11 815.3913
Set-CellValue $ws 34 12 1371.2285
Set-CellValue $ws 34 13 -613.3913
Set-CellValue $ws 34 14 -1775.2285

# Row 132
Set-CellValue $ws 132 8 5953973
Set-CellValue $ws 132 9 1233.8379
Set-CellValue $ws 132 10 17546148
Set-CellValue $ws 132 11 3701.5137
Set-CellValue $ws 132 12 52638444
Set-CellValue $ws 132 13 -1171.5137
Set-CellValue $ws 132 14 -52643504

$ws = $wb.Worksheets.Item("CUL")
# Row 68
Set-CellValue $ws 68 8 2767.0178
Set-CellValue $ws 68 9 797.03125
Set-CellValue $ws 68 10 5393.6665
Set-CellValue $ws 68 11 2391.09375
Set-CellValue $ws 68 12 16180.9995
Set-CellValue $ws 68 13 -1580.09375
Set-CellValue $ws 68 14 -17802.9995

# Row 71
Set-CellValue $ws 71 8 2767.0178
Set-CellValue $ws 71 9 797.03125
Set-CellValue $ws 71 10 5393.6665
Set-CellValue $ws 71 11 7173.28125
Set-CellValue $ws 71 12 48542.9985
Set-CellValue $ws 71 13 -3117.28125
Set-CellValue $ws 71 14 -56654.9985

# Row 123
Set-CellValue $ws 123 8 1400
Set-CellValue $ws 123 9 1400
Set-CellValue $ws 123 11 4200
Set-CellValue $ws 123 13 -1750

# Row 131
Set-CellValue $ws 131 8 8775977
Set-CellValue $ws 131 9 22730338
Set-CellValue $ws 131 10 4570553
Set-CellValue $ws 131 11 68191014
Set-CellValue $ws 131 12 13711659
Set-CellValue $ws 131 13 -68185974
Set-CellValue $ws 131 14 -13721739

$ws = $wb.Worksheets.Item("GSM")
# Row 80
Set-CellValue $ws 80 8 2633.3333
Set-CellValue $ws 80 9 2700
Set-CellValue $ws 80 10 2500
Set-CellValue $ws 80 11 2700
Set-CellValue $ws 80 12 2500
Set-CellValue $ws 80 13 -1702
Set-CellValue $ws 80 14 -4496

# Row 83
Set-CellValue $ws 83 8 2633.3333
Set-CellValue $ws 83 9 2700
Set-CellValue $ws 83 10 2500
Set-CellValue $ws 83 11 13500
Set-CellValue $ws 83 12 12500
Set-CellValue $ws 83 13 -8508
Set-CellValue $ws 83 14 -22484

$ws = $wb.Worksheets.Item("LTW")
# Row 82
Set-CellValue $ws 82 8 1448.75
Set-CellValue $ws 82 9 1100.3334
Set-CellValue $ws 82 10 1797.1666
Set-CellValue $ws 82 11 1100.3334
Set-CellValue $ws 82 12 1797.1666
Set-CellValue $ws 82 13 -739.3334
Set-CellValue $ws 82 14 -2519.1666

# Row 85
Set-CellValue $ws 85 8 1448.75
Set-CellValue $ws 85 9 1100.3334
Set-CellValue $ws 85 10 1797.1666
Set-CellValue $ws 85 11 1100.3334
Set-CellValue $ws 85 12 1797.1666
Set-CellValue $ws 85 13 147.6666
Set-CellValue $ws 85 14 -4293.1666

$ws = $wb.Worksheets.Item("WVR")
# Row 122
Set-CellValue $ws 122 8 20050.26
Set-CellValue $ws 122 9 22693.783
Set-CellValue $ws 122 10 4850
Set-CellValue $ws 122 11 68081.349
Set-CellValue $ws 122 12 14550
Set-CellValue $ws 122 13 -65631.349
Set-CellValue $ws 122 14 -19450

# Row 132
Set-CellValue $ws 132 8 23712.686
Set-CellValue $ws 132 9 40758
Set-CellValue $ws 132 10 5985.56
Set-CellValue $ws 132 12 17956.68
Set-CellValue $ws 132 14 -23016.68
Set-CellValue $ws 132 13 -119744
